$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that must stay literal text (e.g. "416.20", "1.00").
# Assigning such strings straight to .Value lets Excel auto-coerce them to numbers
# (dropping trailing zeros / thousands dots), so we briefly force a Text number format,
# write the value, then reset the style back to Normal so no stray formatting sticks.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.977.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +6.87%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.567.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.85%  "

$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "416.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.38%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.650"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.34%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.556.70"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.73%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.767"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.49%  "

$ws.Range("E11").Value = "  +12.64%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000329"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +44.77%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.84%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.03%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.139.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.91%  "

$ws.Range("E16").Value = "  -0.19%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.39"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.98%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.550.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.28%  "

$ws.Range("E19").Value = "  +5.12%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "66.940.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.80%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.35%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "455.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.75%  "

$ws.Range("E23").Value = "  -2.44%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.22%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.51%  "

$ws.Range("E26").Value = "  +1.42%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.57%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "34.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.93%  "

$ws.Range("E29").Value = "  +1.86%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.79"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.22%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.33"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.61%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.117"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.06%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.32%  "

$ws.Range("E34").Value = "  -4.54%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "40.92"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.24%  "

$ws.Range("E36").Value = "  -0.19%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.53"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.79%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0493"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.50%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0725"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +26.78%  "

$ws.Range("E40").Value = "  +8.78%  "

$ws.Range("E41").Value = "  -0.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.05"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.37%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "148.99"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.37%  "

$ws.Range("E44").Value = "  +0.38%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.46%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.312"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.82%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.39%  "

$ws.Range("E48").Value = "  -4.27%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.34"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.95%  "

$ws.Range("E50").Value = "  +14.88%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "15.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.74%  "
